$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append one new observation record as row 6 (the sheet currently has data
# in rows 1-5; the header is row 1). This mirrors the shape of the existing
# rows: a mushroom-sighting record with taxon info, locality and reporter.

$ws.Range("A6").Value = 112213267
$ws.Range("B6").Value = 90687
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 5964
$ws.Range("F6").Value = "Fjällig taggsvamp s.str."
$ws.Range("G6").Value = "Sarcodon imbricatus s.str."
$ws.Range("H6").Value = "(L.:Fr.) P.Karst."

$ws.Range("P6").Value = "Ol-olssvarttjärnen, Jmt"
$ws.Range("Q6").Value = 446536
$ws.Range("R6").Value = 7032714
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = "Jämtland"
$ws.Range("U6").Value = "Krokom"
$ws.Range("V6").Value = "Jämtland"
$ws.Range("W6").Value = "Alsen"

# Dates are stored as plain text ("yyyy-mm-dd"), not Excel date serials,
# matching how the rest of the sheet encodes Startdatum/Slutdatum.
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-19"
$ws.Range("Y6").Style = "Normal"

$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-19"
$ws.Range("AA6").Style = "Normal"

$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false

$ws.Range("AW6").Value = "Erik Lundmark"
$ws.Range("AX6").Value = "Erik Lundmark"
